$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate a paragraph by (partial) text content rather than a fixed
# index, so the script keeps working even if indices shift.
# ------------------------------------------------------------------
function Get-ParaIndexContainingText([string]$needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1) Remove the two blank paragraphs that sit right above question 2.
# ------------------------------------------------------------------
$q2Index = Get-ParaIndexContainingText("otázka 2")
$emptyIndex = $q2Index - 2

$d.Paragraphs.Item($emptyIndex).Range.Delete()
$d.Paragraphs.Item($emptyIndex).Range.Delete()

# ------------------------------------------------------------------
# 2) Cut out question 2's block: the question paragraph, its three
#    answer paragraphs and the blank paragraph right after them.
# ------------------------------------------------------------------
$q2Index   = Get-ParaIndexContainingText("otázka 2")
$startPara = $d.Paragraphs.Item($q2Index)
$endPara   = $d.Paragraphs.Item($q2Index + 4)   # question + a + b + c + blank

$blockStart = $startPara.Range.Start
$blockEnd   = $endPara.Range.End

$d.Range($blockStart, $blockEnd).Delete()

# ------------------------------------------------------------------
# 3) Re-insert that block right after the (empty) bookmark paragraph
#    that now directly precedes question 3.
# ------------------------------------------------------------------
$q3Index       = Get-ParaIndexContainingText("otázka 3")
$bookmarkIndex = $q3Index - 1
$bookmarkPara  = $d.Paragraphs.Item($bookmarkIndex)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$blockXml = @"
<w:p $wNs><w:pPr><w:pStyle w:val="Odstavecseseznamem"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Testovací</w:t></w:r><w:r><w:t xml:space="preserve"> otázka 2</w:t></w:r></w:p><w:p $wNs><w:pPr><w:pStyle w:val="Odstavecseseznamem"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Testovací</w:t></w:r><w:r><w:t xml:space="preserve"> odpověď a</w:t></w:r></w:p><w:p $wNs><w:pPr><w:pStyle w:val="Odstavecseseznamem"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Testovací</w:t></w:r><w:r><w:t xml:space="preserve"> odpověď b</w:t></w:r></w:p><w:p $wNs><w:pPr><w:pStyle w:val="Odstavecseseznamem"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Testovací</w:t></w:r><w:r><w:t xml:space="preserve"> odpověď c</w:t></w:r></w:p><w:p $wNs/>
"@

$insertionPoint = $d.Range($bookmarkPara.Range.End, $bookmarkPara.Range.End)
[void]$insertionPoint.InsertXML($blockXml)
